# Order components for RevB Micro
# Updates to acquisitions to reflect ordered components (Mouser order placed):
# fills in supplier/part/manufacturer/pricing details for the previously
# placeholder rows, corrects a couple of quantities/values, and appends two
# new LED line items (Red LED / Green LED) that were also ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title block ---
$ws.Range("A2").Value = "Last updated: Oct. 21/14"

# --- Column headers (row 8) ---
$ws.Range("K8").Value = "ROHS"

# --- Row 9: Capacitor 4.7u C203 ---
$ws.Range("A9").Value = "Capacitor"
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = "4.7u"
$ws.Range("D9").Value = "C203"
$ws.Range("E9").Value = "SMD_0603"
$ws.Range("F9").Value = "Mouser"
$ws.Range("G9").Value = "81-GRM18R60J475KE19D"
$ws.Range("H9").Value = "Murata"
$ws.Range("I9").Value = "GRM188R60J475KE19D"
$ws.Range("J9").Value = "0.063/0.035/--"
$ws.Range("K9").Value = "compliant"

# --- Row 10: Microcontroller U101 ---
$ws.Range("A10").Value = "Microcontroller"
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = "N/A"
$ws.Range("D10").Value = "U101"
$ws.Range("E10").Value = "TQFP_32"
$ws.Range("F10").Value = "Mouser"
$ws.Range("G10").Value = "556-ATMEGA328P-AU"
$ws.Range("H10").Value = "Atmel"
$ws.Range("I10").Value = "ATMEGA328P-AU"
$ws.Range("J10").Value = "3.64/2.74/2.58"
$ws.Range("K10").Value = "compliant"

# --- Row 11: Resonator X101 ---
$ws.Range("A11").Value = "Resonator"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "16M"
$ws.Range("D11").Value = "X101"
$ws.Range("E11").Value = "CSTCE"
$ws.Range("F11").Value = "Mouser"
$ws.Range("G11").Value = "81-CSTCE16M0V53-R0"
$ws.Range("H11").Value = "Murata"
$ws.Range("I11").Value = "CSTCE16M0V53-R0"
$ws.Range("J11").Value = "0.434/0.35/--"
$ws.Range("K11").Value = "compliant"

# --- Row 12: Resistor 10k R104 ---
$ws.Range("A12").Value = "Resistor"
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = "10k"
$ws.Range("D12").Value = "R104"
$ws.Range("E12").Value = "SMD_0603"
$ws.Range("F12").Value = "Mouser"
$ws.Range("G12").Value = "71-CRCW0603-2.0K-E3"
$ws.Range("H12").Value = "Vishay"
$ws.Range("I12").Value = "CRCW06032K00FKEA"
$ws.Range("J12").Value = "0.016/0.012/--"
$ws.Range("K12").Value = "exemption"

# --- Row 13: Resistor 180 R106,R107 ---
$ws.Range("A13").Value = "Resistor"
$ws.Range("B13").Value = 14
$ws.Range("C13").Value = 180
$ws.Range("D13").Value = "R106,R107"
$ws.Range("E13").Value = "SMD_0603"
$ws.Range("F13").Value = "Mouser"
$ws.Range("G13").Value = "71-CRCW0603J-180-E3"
$ws.Range("H13").Value = "Vishay"
$ws.Range("I13").Value = "CRCW0603180RJNEA"
$ws.Range("J13").Value = "0.012/0.01/--"
$ws.Range("K13").Value = "exemption"

# --- Row 14: Capacitor 10n C105 ---
$ws.Range("A14").Value = "Capacitor"
$ws.Range("B14").Value = 10
$ws.Range("C14").Value = "10n"
$ws.Range("D14").Value = "C105"
$ws.Range("E14").Value = "SMD_0603"
$ws.Range("F14").Value = "Mouser"
$ws.Range("G14").Value = "77-VJ0603Y103KXAAC"
$ws.Range("H14").Value = "Vishay"
$ws.Range("I14").Value = "VJ0603Y103KXAAC"
$ws.Range("J14").Value = "'0.06/0.036/0.03"
$ws.Range("K14").Value = "compliant"

# --- Row 15: Resistor 0 R105,R108 ---
$ws.Range("A15").Value = "Resistor"
$ws.Range("B15").Value = 15
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = "R105,R108"
$ws.Range("E15").Value = "SMD_0603"
$ws.Range("F15").Value = "Mouser"
$ws.Range("G15").Value = "71-CRCW0603-0-E3"
$ws.Range("H15").Value = "Vishay"
$ws.Range("I15").Value = "CRCW06030000Z0EA"
$ws.Range("J15").Value = "0.012/0.01/--"
$ws.Range("K15").Value = "exemption"

# --- Row 16: Regulator U103 ---
$ws.Range("A16").Value = "Regulator"
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = "3.3V"
$ws.Range("D16").Value = "U103"
$ws.Range("E16").Value = "SOT-23-5"
$ws.Range("F16").Value = "Mouser"
$ws.Range("G16").Value = "595-LP2985-33DBVR"
$ws.Range("H16").Value = "TI"
$ws.Range("I16").Value = "LP2985-33DBVR"
$ws.Range("J16").Value = "0.578/0.43/0.317"
$ws.Range("K16").Value = "compliant"

# --- Row 17: Pin Headers ---
$ws.Range("A17").Value = "Pin Headers"
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = "N/A"
$ws.Range("D17").Value = "N/A"
$ws.Range("E17").Value = "N/A"
$ws.Range("F17").Value = "Mouser"
$ws.Range("G17").Value = "517-9611106404AR"
$ws.Range("H17").Value = "3M"
$ws.Range("I17").Value = "961110-6404-AR"
$ws.Range("J17").Value = "'--/0.437/0.412"
$ws.Range("K17").Value = "compliant"

# --- Row 18: Red LED (replaces old "LEDs"/"!!" placeholder row) ---
$ws.Range("A18").Value = "Red LED"
$ws.Range("B18").Value = 10
$ws.Range("C18").Value = "N/A"
$ws.Range("D18").Value = "LED5"
$ws.Range("E18").Value = "SMD_0603"
$ws.Range("F18").Value = "Mouser"
$ws.Range("G18").Value = "710-150060RS75000"
$ws.Range("H18").Value = "Wurth"
$ws.Range("I18").Value = "150060RS75000"
$ws.Range("J18").Value = "0.253/0.232/--"
$ws.Range("K18").Value = "compliant"

# --- Row 19 (new): Green LED ---
$ws.Range("A19").Value = "Green LED"
$ws.Range("B19").Value = 10
$ws.Range("C19").Value = "N/A"
$ws.Range("D19").Value = "LED4"
$ws.Range("E19").Value = "SMD_0603"
$ws.Range("F19").Value = "Mouser"
$ws.Range("G19").Value = "710-150060VS75000"
$ws.Range("H19").Value = "Wurth"
$ws.Range("I19").Value = "150060VS75000"
$ws.Range("J19").Value = "0.253/0.232/--"
$ws.Range("K19").Value = "compliant"

# --- Restore the view/selection state recorded in the saved workbook ---
$ws.Range("B13").Select()
